# Remove the trailing "Ver no Jupiter..." / "© 2020 ..." footer block
# (and the blank paragraph right before it) that follows the
# "LOM3202: Circuitos Elétricos (Indicação de Conjunto)" line, while
# keeping the blank paragraph that sits right before the final
# page-break paragraph.

$d = $word.ActiveDocument

$anchor = "LOM3202: Circuitos Elétricos (Indicação de Conjunto)"

# Locate the paragraph containing the anchor text.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith($anchor)) {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -gt 0) {
    # The three paragraphs to delete are the blank paragraph, the
    # "Ver no Jupiter..." paragraph and the "© 2020 ..." paragraph,
    # immediately following the anchor paragraph.
    $firstToDelete = $anchorIndex + 1
    $lastToDelete = $anchorIndex + 3

    $start = $d.Paragraphs.Item($firstToDelete).Range.Start
    $end = $d.Paragraphs.Item($lastToDelete).Range.End

    $r = $d.Range($start, $end)
    $r.Delete()
}
